$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing last-header cell's formatting (bold font, thin border,
# centered alignment - style index 1) onto the three new header cells so
# they match the look of the rest of row 1 without inventing a new style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels for the team win/loss/tie record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-48) gets the same 1991 team record: 84 wins, 78
# losses, 0 ties.
$ws.Range("AD2:AD48").Value = 84
$ws.Range("AE2:AE48").Value = 78
$ws.Range("AF2:AF48").Value = 0
